$d = $word.ActiveDocument

# 1. Change the licence text run: "CC BY-NC" -> "CC BY-SA"
$d.Content.Find.Execute("licensed under CC BY-NC 4.0. To view", $true, $false, $false, $false, $false,
                         $true, 1, $false, "licensed under CC BY-SA 4.0. To view", 2)

# 2. Update the hyperlink address + display text for the licence link
foreach ($hl in $d.Hyperlinks) {
    if ($hl.Address -eq "https://creativecommons.org/licenses/by-nc/4.0") {
        $hl.Address = "https://creativecommons.org/licenses/by-sa/4.0"
        $hl.TextToDisplay = "https://creativecommons.org/licenses/by-sa/4.0"
    }
}

# 3. Remove the stray lastRenderedPageBreak before "Follow the instructor directions"
$d.Content.Find.Execute("Follow the instructor directions", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Follow the instructor directions", 2)
